$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.262.87"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.633.63"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.53"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.522"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.33"
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "1.634.59"
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.546"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "65.16"
$ws.Range("E15").Value = "  -3.80%  "
$ws.Range("D16").Value = "27.194.08"
$ws.Range("D17").Value = "0.0₃0737"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.55"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("E21").Value = "  -1.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.46"
$ws.Range("E22").Value = "  -3.93%  "
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.10"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.28"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  -1.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0508"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("D33").Value = "1.317.23"
$ws.Range("E33").Value = "  +3.98%  "
$ws.Range("E34").Value = "  -1.25%  "
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.850"
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.542"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.26"
$ws.Range("E40").Value = "  +2.39%  "
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.99"
$ws.Range("E42").Value = "  +2.92%  "
$ws.Range("D43").Value = "1.770.73"
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.22"
$ws.Range("E44").Value = "  -4.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.72"
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("D47").Value = "0.0₆0107"
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("E48").Value = "  +21.87%  "
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0962"
$ws.Range("E51").Value = "  -0.94%  "
